# Auto-generated script applying the cryptos.xlsx price/volume refresh diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '96.995.82'
$ws.Range("E2").Value = '  +0.26%  '

$ws.Range("D3").Value = '3.702.65'

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.64'
$ws.Range("E5").Value = '  -2.54%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.92'
$ws.Range("E6").Value = '  +2.05%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '656.05'
$ws.Range("E7").Value = '  -2.02%  '

$ws.Range("E8").Value = '  -0.41%  '

$ws.Range("E9").Value = '  -3.12%  '

$ws.Range("E10").Value = '  +0.03%  '

$ws.Range("D11").Value = '3.701.02'
$ws.Range("E11").Value = '  +0.36%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '44.27'
$ws.Range("E12").Value = '  -2.81%  '

$ws.Range("E13").Value = '  +1.28%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000300'
$ws.Range("E14").Value = '  +11.15%  '

$ws.Range("E15").Value = '  +2.21%  '

$ws.Range("D16").Value = '4.393.92'
$ws.Range("E16").Value = '  +0.47%  '

$ws.Range("D17").Value = '96.716.65'
$ws.Range("E17").Value = '  +0.33%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.95'
$ws.Range("E18").Value = '  -1.11%  '

$ws.Range("D19").Value = '3.705.37'
$ws.Range("E19").Value = '  +0.56%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.00'
$ws.Range("E20").Value = '  +1.39%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '18.64'
$ws.Range("E21").Value = '  +1.01%  '

$ws.Range("E22").Value = '  -4.59%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '523.73'
$ws.Range("E23").Value = '  +0.61%  '

$ws.Range("E24").Value = '  -2.02%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000212'
$ws.Range("E25").Value = '  +2.10%  '

$ws.Range("E26").Value = '  -1.56%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '101.96'
$ws.Range("E27").Value = '  -0.91%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.192'
$ws.Range("E28").Value = '  +15.38%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '13.38'
$ws.Range("E29").Value = '  +2.28%  '

$ws.Range("E30").Value = '  +0.90%  '

$ws.Range("E31").Value = '  -1.64%  '

$ws.Range("E32").Value = '  +0.08%  '

$ws.Range("E33").Value = '  +1.07%  '

$ws.Range("E34").Value = '  +1.34%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '32.29'
$ws.Range("E36").Value = '  -1.92%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '647.04'
$ws.Range("E37").Value = '  +3.80%  '

$ws.Range("E38").Value = '  +1.85%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.83'
$ws.Range("E39").Value = '  +0.33%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.84'
$ws.Range("E41").Value = '  +10.24%  '

$ws.Range("E42").Value = '  +4.76%  '

$ws.Range("B43").Value = 'Kaspa'
$ws.Range("C43").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.160'
$ws.Range("E43").Value = '  +0.08%  '

$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '40.43'
$ws.Range("E44").Value = '  -5.44%  '

$ws.Range("E45").Value = '  -0.21%  '

$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0463'
$ws.Range("E46").Value = '  +1.10%  '

$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.450'
$ws.Range("E47").Value = '  +4.33%  '

$ws.Range("E48").Value = '  -0.86%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '23.64'
$ws.Range("E49").Value = '  +0.08%  '

$ws.Range("E50").Value = '  -0.67%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.57'
$ws.Range("E51").Value = '  +1.41%  '
